$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Cypher query text for the "StudyFilesTab" row (literal here-string: no
# interpolation, so backticks / $ / quotes survive untouched).
$studyFilesQuery = @'
 MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['COTC022'] and demo.breed in ['Anatolian Shepherd Dog','Saint Bernard'] and diag.disease_term in ['Osteosarcoma'] and diag.primary_disease_site in ['Bone (Appendicular)']
WITH DISTINCT f, s
RETURN 
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS `File Format`,
  coalesce(f.file_size, '') AS `Size`,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Reuse the exact text already present for the other rows (StatQuery /
# Neo4j-file / Web-file columns are identical across every tab row) so the
# new row shares the same string values as rows 2-4.
$statQuery = $ws.Range("C4").Value2
$neo4jFile = $ws.Range("D4").Value2
$webFile   = $ws.Range("E4").Value2

$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $studyFilesQuery
$ws.Range("C5").Value = $statQuery
$ws.Range("D5").Value = $neo4jFile
$ws.Range("E5").Value = $webFile

# Match the wrap-text style used by column B/C on the existing rows.
$ws.Range("B5:C5").WrapText = $true

# The source row's wrapped text auto-sized to this height in the author's Excel.
$ws.Rows.Item(5).RowHeight = 232

# New active selection lands on the newly added row, same as the source edit.
$ws.Range("C5").Select() | Out-Null
